$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.535.62"
Set-TextValue $ws.Range("E2") "  +0.13%  "
Set-TextValue $ws.Range("D3") "2.492.53"
Set-TextValue $ws.Range("E3") "  -1.21%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "569.94"
Set-TextValue $ws.Range("E5") "  -0.63%  "
Set-TextValue $ws.Range("D6") "166.30"
Set-TextValue $ws.Range("E6") "  -0.42%  "
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("E8") "  -1.43%  "
Set-TextValue $ws.Range("E9") "  -1.02%  "
Set-TextValue $ws.Range("E10") "  -0.69%  "
Set-TextValue $ws.Range("D11") "0.351"
Set-TextValue $ws.Range("E11") "  -1.45%  "
Set-TextValue $ws.Range("D12") "4.87"
Set-TextValue $ws.Range("E12") "  -0.95%  "
Set-TextValue $ws.Range("D13") "2.949.55"
Set-TextValue $ws.Range("E13") "  -1.07%  "
Set-TextValue $ws.Range("D14") "69.359.14"
Set-TextValue $ws.Range("E14") "  +0.20%  "
Set-TextValue $ws.Range("E15") "  -0.83%  "
Set-TextValue $ws.Range("D16") "24.22"
Set-TextValue $ws.Range("E16") "  -2.82%  "
Set-TextValue $ws.Range("D17") "2.505.41"
Set-TextValue $ws.Range("E17") "  -0.64%  "
Set-TextValue $ws.Range("E18") "  -1.10%  "
Set-TextValue $ws.Range("D19") "354.55"
Set-TextValue $ws.Range("E19") "  +1.45%  "
Set-TextValue $ws.Range("D20") "7.39"
Set-TextValue $ws.Range("E20") "  -3.21%  "
Set-TextValue $ws.Range("D21") "3.89"
Set-TextValue $ws.Range("E21") "  -0.82%  "
Set-TextValue $ws.Range("D22") "1.89"
Set-TextValue $ws.Range("E22") "  -5.32%  "
Set-TextValue $ws.Range("E23") "  -0.02%  "
Set-TextValue $ws.Range("E24") "  -1.27%  "
Set-TextValue $ws.Range("D25") "3.80"
Set-TextValue $ws.Range("E25") "  -4.39%  "
Set-TextValue $ws.Range("D27") "8.60"
Set-TextValue $ws.Range("E27") "  -3.76%  "
Set-TextValue $ws.Range("D28") "0.993"
Set-TextValue $ws.Range("E28") "  -0.35%  "
Set-TextValue $ws.Range("D29") "0.0₃0873"
Set-TextValue $ws.Range("E29") "  -2.75%  "
Set-TextValue $ws.Range("D30") "7.61"
Set-TextValue $ws.Range("E30") "  -3.04%  "
Set-TextValue $ws.Range("E31") "  -3.99%  "
Set-TextValue $ws.Range("D32") "437.59"
Set-TextValue $ws.Range("E32") "  -5.71%  "
Set-TextValue $ws.Range("D33") "1.00"
Set-TextValue $ws.Range("E33") "  +0.05%  "
Set-TextValue $ws.Range("E34") "  -1.28%  "
Set-TextValue $ws.Range("D35") "155.01"
Set-TextValue $ws.Range("E35") "  -1.58%  "
Set-TextValue $ws.Range("D36") "0.112"
Set-TextValue $ws.Range("E36") "  -3.60%  "
Set-TextValue $ws.Range("E37") "  +0.26%  "
Set-TextValue $ws.Range("D38") "18.18"
Set-TextValue $ws.Range("E38") "  -2.17%  "
Set-TextValue $ws.Range("E39") "  +0.03%  "
Set-TextValue $ws.Range("D40") "2.62"
Set-TextValue $ws.Range("E40") "  +64.49%  "
Set-TextValue $ws.Range("E41") "  -1.87%  "
Set-TextValue $ws.Range("D42") "4.59"
Set-TextValue $ws.Range("E42") "  -3.07%  "
Set-TextValue $ws.Range("E43") "  -2.11%  "
Set-TextValue $ws.Range("E44") "  -4.29%  "
Set-TextValue $ws.Range("E45") "  -5.64%  "
Set-TextValue $ws.Range("D46") "138.47"
Set-TextValue $ws.Range("E46") "  -2.42%  "
Set-TextValue $ws.Range("D47") "3.42"
Set-TextValue $ws.Range("E47") "  -1.69%  "
Set-TextValue $ws.Range("E48") "  -3.55%  "
Set-TextValue $ws.Range("E49") "  -1.16%  "
Set-TextValue $ws.Range("D50") "0.571"
Set-TextValue $ws.Range("E50") "  -1.59%  "
Set-TextValue $ws.Range("D51") "0.0925"
Set-TextValue $ws.Range("E51") "  -1.08%  "
